$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: type row
$ws.Range("A2").Value = "iaest-measure:regimen-2-digitos"
$ws.Range("J2").Value = "iaest-measure:direccion-provincial-nombre"

# Row 3: dim/medida row
$ws.Range("A3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: class row
$ws.Range("A4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: mapping file row - A5 (mapping-regimen-2-digitos.xlsx) is removed entirely
$ws.Range("A5").Clear()
